$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 98 values
$ws.Cells.Item(98, 1).Value = 24
$ws.Cells.Item(98, 2).Value = 5
$ws.Cells.Item(98, 3).Value = 5

# Update row 99 values
$ws.Cells.Item(99, 1).Value = 21
$ws.Cells.Item(99, 2).Value = 5
$ws.Cells.Item(99, 3).Value = 5

# Delete rows 100 and 101 entirely (shift cells up)
$ws.Range("A100:C101").Delete()

# Update the sheet view (scroll position + active selection)
$ws.Range("E106").Select()
$excel.ActiveWindow.ScrollRow = 91
